$wb = $excel.ActiveWorkbook

# Sheet 1: "VENTAS POR GRUPO"
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("D10").Value = 354.43
$ws1.Range("D22").Value = "1 de 20"

# Sheet 2: "VENTA MENSUAL"
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F10").Value = 354.43
$ws2.Range("F22").Value = 1056.14

# Sheet 3: "CUMPLIMIENTO MENSUAL"
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D2").Value = 71723.53999999999
$ws3.Range("E2").Value = -71723.53999999999
$ws3.Range("D4").Value = 74961.42999999999
$ws3.Range("E4").Value = -59489.87069999999
$ws3.Range("F4").Value = 4.845111507280329
